$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "iaest-measure:estado-de-la-informacion"
$ws.Range("M2").Value = "sdmx-dimension:refArea"

$ws.Range("K3").Value = "medida"
$ws.Range("M3").Value = "dim"

$ws.Range("K4").Value = "xsd:int"
$ws.Range("M4").Value = "URI-Municipio"

$ws.Range("K5").Clear()
